# Apply the edits described by the diff:
# 1. In sheet "classes", insert a new row at A2 with value "speed-bumps-and-potholes",
#    shifting the existing Pothole/RoadPath/SpeedBump rows down by one.
# 2. Update the selection (active cell) in sheet "classes" to D34.
# 3. Update the selection (active cell) in sheet "rules" to B11.

$wb = $excel.ActiveWorkbook

$classes = $wb.Worksheets.Item("classes")
$classes.Activate()
$classes.Rows.Item(2).Insert()
$classes.Range("A2").Value = "speed-bumps-and-potholes"
$classes.Range("D34").Select()

$rules = $wb.Worksheets.Item("rules")
$rules.Activate()
$rules.Range("B11").Select()
